$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Target widths (2.140625 / 3.140625 / 5.7109375) are pre-existing widths used
# elsewhere on this sheet; the closest values reachable through the
# ColumnWidth automation property (which quantizes to a 1/6-character grid)
# are used here.
$ws.Columns.Item(4).ColumnWidth = 2.333333333333333
$ws.Columns.Item(5).ColumnWidth = 1.333333333333333
$ws.Columns.Item(7).ColumnWidth = 1.333333333333333
$ws.Columns.Item(17).ColumnWidth = 4.833333333333334

# --- Row 1 cell value changes ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 31
$ws.Range("D1").Value = 20
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 16
$ws.Range("G1").Value = 9
$ws.Range("H1").Value = 15
$ws.Range("I1").Value = 20
$ws.Range("J1").Value = 32
$ws.Range("K1").Value = 27
$ws.Range("L1").Value = 33
$ws.Range("M1").Value = 0.033
$ws.Range("N1").Value = 0.048
$ws.Range("O1").Value = 0.015
$ws.Range("P1").Value = 0.075
$ws.Range("Q1").Value = 0.06 - 0.039
